$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Header row (row 1): Vietnamese template headers, reordered ---
$ws1.Range("A1").Value = "Mã(*)"
$ws1.Range("B1").Value = "Họ và Tên(*)"
$ws1.Range("C1").Value = "Email(*)"
$ws1.Range("D1").Value = "Giới tính"

# --- Sample data row (row 2) ---
$ws1.Range("A2").NumberFormat = "@"
$ws1.Range("A2").Value = "19110001"
$ws1.Range("B2").Value = "Nguyen Van A"
$ws1.Range("C2").Value = "19110001@student.hcmute.edu.vn"
$ws1.Range("D2").Value = "male"

# --- Re-point the e-mail hyperlink at the new sample address ---
$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("C2"), "mailto:19110001@student.hcmute.edu.vn")

# --- Bump the whole used area to a 12pt font (matches the new template look) ---
$ws1.Range("A1:D2").Font.Size = 12
$ws1.Range("B4:D4").Font.Size = 12

# --- Clear the stray A4 formatting left over from the old template ---
$ws1.Range("A4").Clear()

# --- Column widths for the new, wider Vietnamese labels ---
$ws1.Columns.Item(1).ColumnWidth = 21.67
$ws1.Columns.Item(2).ColumnWidth = 51.17
$ws1.Columns.Item(3).ColumnWidth = 54.5
$ws1.Columns.Item(4).ColumnWidth = 23.33
